# Add I0 (column I) and IF (column J) data, matching the commit
# "I0 and IF added".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the format from H1 (bold/border/centered) onto I1:J1,
# then set the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..66: (row, I value, J value)
$data = @(
    @(2, 9, 10),
    @(3, 9, 9),
    @(4, 6, 6),
    @(5, 8, 8),
    @(6, 9, 9),
    @(7, 7, 7),
    @(8, 6, 6),
    @(9, 8, 8),
    @(10, 7, 8),
    @(11, 6, 6),
    @(12, 8, 9),
    @(13, 5, 6),
    @(14, 7, 8),
    @(15, 6, 7),
    @(16, 3, 6),
    @(17, 9, 9),
    @(18, 9, 9),
    @(19, 5, 6),
    @(20, 9, 9),
    @(21, 9, 9),
    @(22, 6, 7),
    @(23, 9, 9),
    @(24, 6, 6),
    @(25, 1, 3),
    @(26, 5, 6),
    @(27, 7, 7),
    @(28, 8, 8),
    @(29, 9, 9),
    @(30, 8, 8),
    @(31, 7, 8),
    @(32, 7, 8),
    @(33, 10, 10),
    @(34, 8, 8),
    @(35, 7, 7),
    @(36, 6, 7),
    @(37, 7, 7),
    @(38, 8, 9),
    @(39, 8, 8),
    @(40, 9, 9),
    @(41, 7, 8),
    @(42, 8, 8),
    @(43, 9, 10),
    @(44, 7, 7),
    @(45, 8, 9),
    @(46, 9, 9),
    @(47, 8, 8),
    @(48, 8, 8),
    @(49, 7, 7),
    @(50, 8, 8),
    @(51, 7, 7),
    @(52, 9, 9),
    @(53, 9, 10),
    @(54, 9, 9),
    @(55, 8, 8),
    @(56, 9, 9),
    @(57, 9, 9),
    @(58, 8, 8),
    @(59, 9, 9),
    @(60, 8, 8),
    @(61, 7, 7),
    @(62, 7, 7),
    @(63, 5, 5),
    @(64, 6, 6),
    @(65, 7, 7),
    @(66, 3, 3)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Output "I0/IF columns populated"
